$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking text values (e.g. '1.005', '26.442.51') that must stay
# exact strings matching the source feed's formatting. Force Text number format on these
# cells before assigning the value so Excel/COM does not silently coerce them into real
# floating point numbers (which would drop trailing zeros, switch to scientific notation
# for very small numbers, etc). Values with two or more dots (e.g. '26.442.51') are not
# parseable as numbers anyway, but we format all Column D cells uniformly for consistency.
$dCells = @('D2','D3','D4','D5','D6','D7','D8','D9','D10','D12','D13','D14','D15','D16','D17','D18','D19','D20','D21','D22','D23','D24','D25','D26','D27','D28','D29','D30','D31','D32','D33','D34','D35','D36','D37','D38','D39','D40','D41','D42','D43','D44','D45','D46','D47','D48','D49','D50','D51')
foreach ($ref in $dCells) {
    $ws.Range($ref).NumberFormat = '@'
}

# Column D (Price) updates
$ws.Range('D2').Value = '26.442.51'
$ws.Range('D3').Value = '1.687.23'
$ws.Range('D4').Value = '1.005'
$ws.Range('D5').Value = '219.33'
$ws.Range('D6').Value = '0.5065'
$ws.Range('D7').Value = '1.005'
$ws.Range('D8').Value = '0.2668'
$ws.Range('D9').Value = '22.03'
$ws.Range('D10').Value = '0.06294'
$ws.Range('D12').Value = '1.689.57'
$ws.Range('D13').Value = '4.542'
$ws.Range('D14').Value = '0.5777'
$ws.Range('D15').Value = '1.917.46'
$ws.Range('D16').Value = '0.000008584'
$ws.Range('D17').Value = '65.18'
$ws.Range('D18').Value = '26.498.29'
$ws.Range('D19').Value = '5.001'
$ws.Range('D20').Value = '1.005'
$ws.Range('D21').Value = '10.89'
$ws.Range('D22').Value = '185.66'
$ws.Range('D23').Value = '6.258'
$ws.Range('D24').Value = '1.006'
$ws.Range('D25').Value = '144.65'
$ws.Range('D26').Value = '7.491'
$ws.Range('D27').Value = '0.1164'
$ws.Range('D28').Value = '15.68'
$ws.Range('D29').Value = '1.341'
$ws.Range('D30').Value = '0.05741'
$ws.Range('D31').Value = '1.334'
$ws.Range('D32').Value = '3.526'
$ws.Range('D33').Value = '3.518'
$ws.Range('D34').Value = '1.660'
$ws.Range('D35').Value = '1.016'
$ws.Range('D36').Value = '0.5950'
$ws.Range('D37').Value = '2.359'
$ws.Range('D38').Value = '2.676'
$ws.Range('D39').Value = '1.102.27'
$ws.Range('D40').Value = '0.01610'
$ws.Range('D41').Value = '0.8651'
$ws.Range('D42').Value = '5.898'
$ws.Range('D43').Value = '1.004'
$ws.Range('D44').Value = '99.89'
$ws.Range('D45').Value = '1.845.38'
$ws.Range('D46').Value = '0.00000000117'
$ws.Range('D47').Value = '56.31'
$ws.Range('D48').Value = '1.004'
$ws.Range('D49').Value = '8.034'
$ws.Range('D50').Value = '0.4315'
$ws.Range('D51').Value = '0.05216'

# Row 39/40 Coin name + Link swap (VeChain/Maker order changed)
$ws.Range('B39').Value = 'Maker'
$ws.Range('C39').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'

# Column E (Volume 1h %) updates - these already contain leading/trailing spaces and a
# '%' sign so Excel keeps them as plain text without any extra formatting needed.
$ws.Range('E2').Value = '  -7.50%  '
$ws.Range('E3').Value = '  -5.98%  '
$ws.Range('E4').Value = '  +0.28%  '
$ws.Range('E5').Value = '  -5.28%  '
$ws.Range('E6').Value = '  -14.14%  '
$ws.Range('E7').Value = '  +0.14%  '
$ws.Range('E8').Value = '  -3.89%  '
$ws.Range('E9').Value = '  -5.75%  '
$ws.Range('E10').Value = '  -6.97%  '
$ws.Range('E11').Value = '  -2.44%  '
$ws.Range('E12').Value = '  -5.78%  '
$ws.Range('E13').Value = '  -5.45%  '
$ws.Range('E14').Value = '  -5.87%  '
$ws.Range('E15').Value = '  -5.86%  '
$ws.Range('E16').Value = '  -3.86%  '
$ws.Range('E17').Value = '  -13.93%  '
$ws.Range('E18').Value = '  -7.20%  '
$ws.Range('E19').Value = '  -7.83%  '
$ws.Range('E20').Value = '  +0.16%  '
$ws.Range('E21').Value = '  -5.21%  '
$ws.Range('E22').Value = '  -11.23%  '
$ws.Range('E23').Value = '  -8.48%  '
$ws.Range('E24').Value = '  +0.10%  '
$ws.Range('E25').Value = '  -5.12%  '
$ws.Range('E26').Value = '  -6.89%  '
$ws.Range('E27').Value = '  -7.89%  '
$ws.Range('E28').Value = '  -4.63%  '
$ws.Range('E29').Value = '  -5.06%  '
$ws.Range('E30').Value = '  -6.87%  '
$ws.Range('E31').Value = '  -6.22%  '
$ws.Range('E32').Value = '  -6.98%  '
$ws.Range('E33').Value = '  -6.63%  '
$ws.Range('E34').Value = '  -4.06%  '
$ws.Range('E35').Value = '  -3.08%  '
$ws.Range('E36').Value = '  -7.32%  '
$ws.Range('E37').Value = '  -5.76%  '
$ws.Range('E38').Value = '  -1.19%  '
$ws.Range('E39').Value = '  -3.91%  '
$ws.Range('E40').Value = '  -4.94%  '
$ws.Range('E41').Value = '  -0.95%  '
$ws.Range('E42').Value = '  -6.87%  '
$ws.Range('E43').Value = '  +0.01%  '
$ws.Range('E44').Value = '  -0.56%  '
$ws.Range('E45').Value = '  -5.27%  '
$ws.Range('E46').Value = '  +7.13%  '
$ws.Range('E47').Value = '  -6.49%  '
$ws.Range('E48').Value = '  +0.16%  '
$ws.Range('E49').Value = '  -3.70%  '
$ws.Range('E50').Value = '  -3.58%  '
$ws.Range('E51').Value = '  -4.29%  '
